# data : case 1
# Replaces the A1:B32 sample-data block with a new set of values and
# narrows column B from 15.42578125 to ~14.71 (OOXML width units).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    ,@(-0.11449925658460813, 0.11401871104348515)
    ,@(-0.028031130590594699, 0.026897829166939502)
    ,@(0.1047432272631319, -0.10528701869039736)
    ,@(-0.1787058612567094, 0.17768106066594669)
    ,@(-0.17168106106860126, 0.16961622171553703)
    ,@(-0.086515045214242914, 0.0864009664743568)
    ,@(-0.066400966969618835, 0.066146086033311846)
    ,@(-0.0461460865330654, 0.045966618346450439)
    ,@(-0.039966618774903928, 0.039826610767492276)
    ,@(-0.033826611200218792, 0.03381103779799588)
    ,@(-0.051480473114548175, 0.051400889730288668)
    ,@(-0.045400890165177898, 0.045153316195359672)
    ,@(-0.039153316637627888, 0.039085945307554937)
    ,@(-0.027085945784699028, 0.027053452295812797)
    ,@(-0.021053452741948142, 0.02102788373620168)
    ,@(-0.015027884183895779, 0.015004434200772909)
    ,@(-0.0090044346505218087, 0.0089999995333851501)
    ,@(-0.036109857873483975, 0.0360966389975097)
    ,@(-0.027096639411961387, 0.027013638983043098)
    ,@(-0.018013639401175396, 0.018004262968839058)
    ,@(-0.0090042633875508926, 0.0089999995808653921)
    ,@(-0.12264105684582027, 0.12206459378089463)
    ,@(-0.084628630472859179, 0.084125641182657418)
    ,@(-0.042125641789492185, 0.041999999389831011)
    ,@(-0.07710117691864582, 0.076987441122135181)
    ,@(-0.026184836664967293, 0.026151293177839818)
    ,@(-0.020151293597007403, 0.020048225487221849)
    ,@(-0.014048225908595668, 0.013987819330056439)
    ,@(-0.0019878197852296609, 0.0019699915811060009)
    ,@(0.018030007920097457, -0.018136296179442457)
    ,@(-0.027018745430941848, 0.027000766303313384)
    ,@(-0.0060007668065598097, 0.0059999995782975546)
)

for ($i = 0; $i -lt $newValues.Count; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $newValues[$i][0]
    $ws.Cells.Item($r, 2).Value = $newValues[$i][1]
}

# Column B width: 14 characters (renders as 14.7109375 in OOXML units).
$ws.Columns.Item(2).ColumnWidth = 13.83
